$d = $word.ActiveDocument

# Locate the trailing empty "List Paragraph"-styled paragraph that currently
# holds the _GoBack bookmark (last empty ListParagraph item in the To-do list).
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text
    if ($styleName -eq "List Paragraph" -and $text.Trim().Length -eq 0) {
        $targetIndex = $i
    }
}

# Use the preceding (already-numbered) list item as the numbering template so
# new items continue the same list (numId) instead of starting a new one.
$templatePara = $d.Paragraphs.Item($targetIndex - 1)
$template = $templatePara.Range.ListFormat.ListTemplate

$target = $d.Paragraphs.Item($targetIndex)

# 1) Insert a brand-new empty "List Paragraph" after the target - it stays
#    plain (no numbering), matching the trailing blank item in the diff.
$target.Range.InsertParagraphAfter()

# 2) Insert a brand-new paragraph before the target for "Logout fix" and give
#    it the list numbering.
$target.Range.InsertParagraphBefore()
$logoutPara = $d.Paragraphs.Item($targetIndex)
$logoutPara.Range.Text = "Logout fix"
$logoutPara.Range.ListFormat.ApplyListTemplate($template, $true)

# 3) The original bookmark-holding paragraph is now shifted by one; give it
#    the "Input boxes..." text and matching numbering, keeping the bookmark
#    that already lives in that paragraph intact.
$target = $d.Paragraphs.Item($targetIndex + 1)
$target.Range.InsertBefore("Input boxes need to be equal except for message one")
$target.Range.ListFormat.ApplyListTemplate($template, $true)
